# Actualizacoes dia 27 as 17
# Re-sort the CHIMANIMANI institution rows (2-8) into their new order while
# keeping each institution's own Masculino/Feminino/TOTAL values intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target order (row -> Instituicao, Masculino, Feminino, TOTAL)
$rows = @(
    @{ Row = 2; Inst = "ISPM";                   M = 0;  F = 1;  T = 1  },
    @{ Row = 3; Inst = "SDAE SUSSUNDENGA";        M = 1;  F = 0;  T = 1  },
    @{ Row = 4; Inst = "MICAIA";                  M = 3;  F = 3;  T = 6  },
    @{ Row = 5; Inst = "ITAM";                    M = 1;  F = 0;  T = 1  },
    @{ Row = 6; Inst = "PARQUE DE CHIMANIMANI";   M = 1;  F = 0;  T = 1  },
    @{ Row = 7; Inst = "UCM";                     M = 1;  F = 0;  T = 1  },
    @{ Row = 8; Inst = "UNIZAMBEZE";              M = 0;  F = 1;  T = 1  }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Inst
    $ws.Cells.Item($r.Row, 3).Value = $r.M
    $ws.Cells.Item($r.Row, 4).Value = $r.F
    $ws.Cells.Item($r.Row, 5).Value = $r.T
}

$wb.Save()
